# Adapt column header formatting to respective input file names (#7)
#   "<header>_old" -> "<header>_FV2310"
#   "<header>_new" -> "<header>_FV2404"
# Also freeze the header row and turn the data range into a proper Excel
# Table ("Table1") with an AutoFilter, matching the regenerated AHB export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Freeze the header row (pane split after row 1).
$sel = $ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true

# Turn the used range into a real Excel Table so the header row gets a
# filter and the table shows up as "Table1" in the workbook.
$rng = $ws.Range("A1:U74")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

Write-Host "Applied FV2310/FV2404 header rename, froze header row, added Table1."
